$d = $word.ActiveDocument

# Locate the two paragraphs that should be removed:
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#    pages. Original theme under Creative Commons Attribution"
# They immediately follow the "MÁQUINASEscola PRO-TEC" paragraph (and a
# blank paragraph after it), and are themselves followed by another blank
# paragraph before the trailing page-break paragraph.
$firstIdx = -1
$lastIdx = -1
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $firstIdx = $i
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $lastIdx = $i
    }
}

if ($firstIdx -gt 0 -and $lastIdx -ge $firstIdx) {
    # Deleting a paragraph's full Range (text + paragraph mark) merges it
    # away and shifts subsequent paragraphs down by one index, so
    # repeatedly deleting at $firstIdx removes the whole run.
    $count = $lastIdx - $firstIdx + 1
    for ($k = 0; $k -lt $count; $k++) {
        $p = $d.Paragraphs.Item($firstIdx)
        $p.Range.Delete()
    }
}
